$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(-16.8282976673956, 2.524515210142058, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956),
  @(-16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, 2.700494077645235, -16.8282976673956, -16.8282976673956),
  @(-16.8282976673956, 2.566058258069836, -16.8282976673956, -16.8282976673956, 3.346965571009611, -16.8282976673956, 1.620500226552398, -16.8282976673956, 2.364045721435527, -16.8282976673956),
  @(-16.8282976673956, 2.08754962474664, -16.8282976673956, -16.8282976673956, -16.8282976673956, 2.890437538351021, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956),
  @(-16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956),
  @(2.689927093857762, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956),
  @(-16.8282976673956, -16.8282976673956, -16.8282976673956, 1.852491067669088, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956),
  @(3.759905471700022, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956),
  @(-16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, 1.323952509562591, -16.8282976673956, 1.896845548839924),
  @(-16.8282976673956, -16.8282976673956, -16.8282976673956, 2.796877842060882, -16.8282976673956, 2.702704868371714, -16.8282976673956, -16.8282976673956, -16.8282976673956, 1.707421951278346),
  @(-16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956),
  @(-16.8282976673956, -16.8282976673956, -16.8282976673956, 2.463468844547994, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, 2.050917610727752, 1.884877373893706),
  @(-16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, 2.264296901071354),
  @(-16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956),
  @(-16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, 2.230721934980828, -16.8282976673956),
  @(-16.8282976673956, 0.8698982880806445, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, 1.092851735128878, 1.800905959923984, 1.767502500335096, -16.8282976673956),
  @(-16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, 1.411659738320312, 1.077776719950764, 1.383899455213165, -16.8282976673956),
  @(-16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, -16.8282976673956, 1.622434345056547, 1.645285220907652, -16.8282976673956, -16.8282976673956),
  @(-16.8282976673956, 0.1460384934540512, 4.321916316724827, -16.8282976673956, 3.296425702174786, -16.8282976673956, 2.001242365945034, 1.184506891609822, -16.8282976673956, 2.173992698365443),
  @(-16.8282976673956, 0.1888695957185318, -16.8282976673956, 1.972276826247848, -16.8282976673956, 2.602810171779047, 2.336081445661948, -16.8282976673956, -16.8282976673956, -16.8282976673956)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowValues[$c]
    }
}
